$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1594.2106
$ws.Range("J17").Value = 1340.8
$ws.Range("L17").Value = 4022.4
$ws.Range("N17").Value = -4358.4
$ws.Range("H58").Value = 943.5
$ws.Range("I58").Value = 591.6667
$ws.Range("K58").Value = 1775.0001
$ws.Range("M58").Value = -1625.0001
$ws.Range("H98").Value = 3162.85
$ws.Range("I98").Value = 3123.7334
$ws.Range("K98").Value = 3123.7334
$ws.Range("M98").Value = -1625.7334
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("N115").Value = 0
$ws.Range("H122").Value = 3162.85
$ws.Range("I122").Value = 3123.7334
$ws.Range("K122").Value = 9371.200199999999
$ws.Range("M122").Value = -6921.200199999999
$ws.Range("H129").Value = 858.9231
$ws.Range("J129").Value = 971.0526
$ws.Range("L129").Value = 2913.1578
$ws.Range("N129").Value = -12913.1578
$ws.Range("H135").Value = 448.5
$ws.Range("I135").Value = 410.90625
$ws.Range("K135").Value = 3698.15625
$ws.Range("M135").Value = -1163.15625
$ws.Range("H137").Value = 1439.6471
$ws.Range("I137").Value = 960.5454999999999
$ws.Range("K137").Value = 2881.6365
$ws.Range("M137").Value = -331.6364999999996
$ws.Range("H141").Value = 596949
$ws.Range("I141").Value = 683566.2
$ws.Range("J141").Value = 5064.8335
$ws.Range("K141").Value = 2050698.6
$ws.Range("L141").Value = 15194.5005
$ws.Range("M141").Value = -2045518.6
$ws.Range("N141").Value = -25554.5005
$ws.Range("L115").ClearContents()
$ws.Range("M115").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3545.5132
$ws.Range("I32").Value = 2846.4775
$ws.Range("K32").Value = 2846.4775
$ws.Range("M32").Value = -2559.4775
$ws.Range("H61").Value = 62501370
$ws.Range("I61").Value = 33334794
$ws.Range("K61").Value = 33334794
$ws.Range("M61").Value = -33334582
$ws.Range("H74").Value = 1156.4286
$ws.Range("I74").Value = 947.36664
$ws.Range("K74").Value = 947.36664
$ws.Range("M74").Value = -73.36663999999996
$ws.Range("H77").Value = 1156.4286
$ws.Range("I77").Value = 947.36664
$ws.Range("K77").Value = 4736.8332
$ws.Range("M77").Value = -368.8332
$ws.Range("H110").Value = 2560.9092
$ws.Range("I110").Value = 1815.7
$ws.Range("J110").Value = 10013
$ws.Range("K110").Value = 1815.7
$ws.Range("L110").Value = 10013
$ws.Range("M110").Value = 229.3
$ws.Range("N110").Value = -14103
$ws.Range("H122").Value = 1071.3784
$ws.Range("I122").Value = 976.2646999999999
$ws.Range("J122").Value = 2149.3333
$ws.Range("K122").Value = 2928.7941
$ws.Range("L122").Value = 6447.999899999999
$ws.Range("M122").Value = -478.7941000000001
$ws.Range("N122").Value = -11347.9999
$ws.Range("H132").Value = 1295.125
$ws.Range("I132").Value = 929.63416
$ws.Range("J132").Value = 3435.8572
$ws.Range("K132").Value = 2788.90248
$ws.Range("L132").Value = 10307.5716
$ws.Range("M132").Value = -258.9024799999997
$ws.Range("N132").Value = -15367.5716
$ws.Range("H136").Value = 62501370
$ws.Range("I136").Value = 33334794
$ws.Range("K136").Value = 100004382
$ws.Range("M136").Value = -100001832

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H134").Value = 5025.1875
$ws.Range("I134").Value = 4213.841
$ws.Range("K134").Value = 12641.523
$ws.Range("M134").Value = -10106.523
$ws.Range("M45").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2382592.5
$ws.Range("I31").Value = 7937775.5
$ws.Range("J31").Value = 1799.762
$ws.Range("K31").Value = 7937775.5
$ws.Range("L31").Value = 1799.762
$ws.Range("M31").Value = -7937480.5
$ws.Range("N31").Value = -2389.762
$ws.Range("H34").Value = 2382592.5
$ws.Range("I34").Value = 7937775.5
$ws.Range("J34").Value = 1799.762
$ws.Range("K34").Value = 7937775.5
$ws.Range("L34").Value = 1799.762
$ws.Range("M34").Value = -7937573.5
$ws.Range("N34").Value = -2203.762
$ws.Range("H53").Value = 55555
$ws.Range("J53").Value = 55555
$ws.Range("L53").Value = 55555
$ws.Range("N53").Value = -56769
$ws.Range("H58").Value = 1146439
$ws.Range("I58").Value = 1891950.9
$ws.Range("J58").Value = 3320.6667
$ws.Range("K58").Value = 1891950.9
$ws.Range("L58").Value = 3320.6667
$ws.Range("M58").Value = -1891747.9
$ws.Range("N58").Value = -3726.6667
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920
$ws.Range("H132").Value = 1367.0588
$ws.Range("I132").Value = 824.5454999999999
$ws.Range("J132").Value = 4777.143
$ws.Range("K132").Value = 2473.6365
$ws.Range("L132").Value = 14331.429
$ws.Range("M132").Value = 56.36350000000039
$ws.Range("N132").Value = -19391.429
$ws.Range("H136").Value = 1146439
$ws.Range("I136").Value = 1891950.9
$ws.Range("J136").Value = 3320.6667
$ws.Range("K136").Value = 5675852.699999999
$ws.Range("L136").Value = 9962.000100000001
$ws.Range("M136").Value = -5673302.699999999
$ws.Range("N136").Value = -15062.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1049.5
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("H107").Value = 257.45
$ws.Range("I107").Value = 308
$ws.Range("J107").Value = 206.9
$ws.Range("K107").Value = 924
$ws.Range("L107").Value = 620.7
$ws.Range("M107").Value = 996
$ws.Range("N107").Value = -4460.7
$ws.Range("H122").Value = 931.8182
$ws.Range("J122").Value = 959.5
$ws.Range("L122").Value = 8635.5
$ws.Range("N122").Value = -13535.5
$ws.Range("H134").Value = 3345.9443
$ws.Range("I134").Value = 1175.5714
$ws.Range("J134").Value = 4727.091
$ws.Range("K134").Value = 3526.7142
$ws.Range("L134").Value = 14181.273
$ws.Range("M134").Value = 1543.2858
$ws.Range("N134").Value = -24321.273
$ws.Range("H137").Value = 5002.846
$ws.Range("J137").Value = 5892.7
$ws.Range("L137").Value = 17678.1
$ws.Range("N137").Value = -27878.1
$ws.Range("H138").Value = 2088.6667
$ws.Range("I138").Value = 1916.4286
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 5749.2858
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = -609.2857999999997
$ws.Range("N138").Value = -23780
$ws.Range("H139").Value = 10885.417
$ws.Range("I139").Value = 12562.5
$ws.Range("J139").Value = 2500
$ws.Range("K139").Value = 37687.5
$ws.Range("L139").Value = 7500
$ws.Range("M139").Value = -32547.5
$ws.Range("N139").Value = -17780
$ws.Range("H141").Value = 4617.6665
$ws.Range("I141").Value = 4617.6665
$ws.Range("K141").Value = 13852.9995
$ws.Range("M141").Value = -8672.999500000002
$ws.Range("N11").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 48.125
$ws.Range("I2").Value = 11.153846
$ws.Range("J2").Value = 91.818184
$ws.Range("K2").Value = 11.153846
$ws.Range("L2").Value = 91.818184
$ws.Range("M2").Value = 101.846154
$ws.Range("N2").Value = -317.818184
$ws.Range("H21").Value = 70007
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H30").Value = 70007
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("H132").Value = 1014570.44
$ws.Range("I132").Value = 1480947.9
$ws.Range("K132").Value = 4442843.699999999
$ws.Range("M132").Value = -4440313.699999999
$ws.Range("M21").ClearContents()
$ws.Range("M30").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2300.818
$ws.Range("I22").Value = 2545.8
$ws.Range("J22").Value = 2096.6667
$ws.Range("K22").Value = 2545.8
$ws.Range("L22").Value = 2096.6667
$ws.Range("M22").Value = -2250.8
$ws.Range("N22").Value = -2686.6667
$ws.Range("H27").Value = 2300.818
$ws.Range("I27").Value = 2545.8
$ws.Range("J27").Value = 2096.6667
$ws.Range("K27").Value = 2545.8
$ws.Range("L27").Value = 2096.6667
$ws.Range("M27").Value = -2438.8
$ws.Range("N27").Value = -2310.6667
$ws.Range("H55").Value = 413.6111
$ws.Range("I55").Value = 365
$ws.Range("K55").Value = 365
$ws.Range("M55").Value = -192
$ws.Range("H61").Value = 2521.7144
$ws.Range("I61").Value = 2108.6667
$ws.Range("K61").Value = 2108.6667
$ws.Range("M61").Value = -1906.6667
$ws.Range("H113").Value = 2521.7144
$ws.Range("I113").Value = 2108.6667
$ws.Range("K113").Value = 2108.6667
$ws.Range("M113").Value = 61.33329999999978
$ws.Range("H132").Value = 1347.8116
$ws.Range("I132").Value = 865.35187
$ws.Range("K132").Value = 2596.05561
$ws.Range("M132").Value = -66.05560999999989

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1999.5
$ws.Range("J107").Value = 1999.5
$ws.Range("L107").Value = 5998.5
$ws.Range("N107").Value = -9838.5
$ws.Range("H113").Value = 1999.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1999.6666
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = 5998.9998
$ws.Range("N113").Value = -10338.9998
$ws.Range("H132").Value = 1181.7721
$ws.Range("I132").Value = 769.2857
$ws.Range("J132").Value = 2805.9375
$ws.Range("K132").Value = 2307.8571
$ws.Range("L132").Value = 8417.8125
$ws.Range("M132").Value = 222.1428999999998
$ws.Range("N132").Value = -13477.8125
$ws.Range("L113").ClearContents()
